$d = $word.ActiveDocument

# 1. Fix spelling: "trophys" -> "trophies"
$d.Content.Find.Execute("trophys", $true, $false, $false, $false, $false, $true, 1, $false, "trophies", 2) | Out-Null

# 2. Append a new closing sentence to the end of the "development choices" paragraph
$d.Content.Find.Execute(
    "texturing in this project. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "texturing in this project. And finally, It had enough sentimental value so I would be excited to work on it.",
    2
) | Out-Null

# 3. Fix typo: "each objects normal" -> "each object's normal"
$d.Content.Find.Execute(
    "each objects normal",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "each object's normal",
    2
) | Out-Null

# 4. Fix spelling: "determined by wether" -> "determined by whether"
$d.Content.Find.Execute(
    "determined by wether",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "determined by whether",
    2
) | Out-Null
